# Automatische test-sync: 2025-06-17 22:42:31
# Appends the new mail-log entry (row 56) to the "Logs" sheet, extends the
# conditional-formatting ranges to cover it, and bumps the "Informatieaanvraag"
# tally on the "Dashboard" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet   # "Logs" is the active sheet in this workbook

$newRow = 56

$ws.Cells.Item($newRow, 1).Value = "Wat zijn jullie openingstijden?"
$ws.Cells.Item($newRow, 2).Value = "mailmind.test@zohomail.eu"
$ws.Cells.Item($newRow, 3).Value = "Hallo, ik zou graag willen weten wat jullie openingstijden zijn. Dank je wel!"
$ws.Cells.Item($newRow, 4).Value = "Informatieaanvraag"
$ws.Cells.Item($newRow, 5).Value = "Beste klant,
Bedankt voor uw e-mail. Onze openingstijden zijn van maandag tot en met vrijdag van 9:00 tot 18:00 uur. Op zaterdag zijn wij geopend van 10:00 tot 17:00 uur. Op zondag zijn wij gesloten. Mocht u nog verdere vragen hebben, dan hoor ik dat graag.
Met vriendelijke groet,
[Naam van de assistent]"
$ws.Cells.Item($newRow, 6).Value = "2025-06-17 22:41:31"
$ws.Cells.Item($newRow, 7).Value = "Ja"

# The multi-line answer in column E would otherwise trigger an automatic
# row-height bump (wrap-based re-measure); AutoFit puts row 56 back in line
# with every other (unstyled, default-height) row in the sheet.
$ws.Rows.Item($newRow).AutoFit()

# Grow both conditional-formatting blocks (Categorie / Beantwoord columns) so
# they keep covering the full data range, D2:D56 and G2:G56.
$ws.Range("D2:D55").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("D2:D56"))
$ws.Range("G2:G55").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("G2:G56"))

# Dashboard: "Informatieaanvraag" count goes from 20 to 21 with the new row.
$wsDash = $wb.Worksheets.Item("Dashboard")
$wsDash.Cells.Item(2, 2).Value = 21
